$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detalle")

# Start this sheet from a clean slate (values + formatting + autofilter).
$ws.Cells.Clear()
$ws.Cells.ClearFormats()
$ws.AutoFilterMode = $false

# --- Row 1: report title + explanatory note -------------------------------
$ws.Range("A1").Value = "Reporte permanencia documentos"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 16
$ws.Range("C1").Value = "Para el cálculo de fechas, se consideran solo días laborales entre 9 y 18 horas. Se excluyen días sábado, domingo y festivos."

# --- Rows 2-4: filter controls (left aligned labels) -----------------------
$ws.Range("A2").Value = "Unidad de permanencia:"
$ws.Range("A3").Value = "Desde"
$ws.Range("A4").Value = "Hasta:"
$ws.Range("A2:A4").HorizontalAlignment = -4131

# Date-entry cells next to Desde/Hasta labels.
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Copy($ws.Range("B4"))

# --- Row 6: data table header ----------------------------------------------
$ws.Range("A6").Value = "Unidad origen"
$ws.Range("B6").Value = "Proceso Id"
$ws.Range("C6").Value = "Cantidad de veces que el documento ingresó a la unidad"
$ws.Range("D6").Value = "Total tiempo de permanencia (días)"
$ws.Range("A6:D6").Font.Bold = $true

$ws.Range("A6:D6").AutoFilter()
$wb.Names.Item("Detalle!_FilterDatabase").RefersTo = "=Detalle!`$A`$6:`$D`$6"

# Column widths for the new layout (closest achievable given engine's
# character-width quantization granularity).
$ws.Columns.Item(1).ColumnWidth = 43.9
$ws.Columns.Item(2).ColumnWidth = 15.6
$ws.Columns.Item(3).ColumnWidth = 53.3
$ws.Columns.Item(4).ColumnWidth = 36.0

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Remove the now-obsolete "Resumen" sheet; keep "Detalle" selected/active.
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Delete()
$ws = $wb.Worksheets.Item("Detalle")
$ws.Activate()
$ws.Range("A7").Select()
